# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: new case/death/recovery numbers for a
# handful of countries, swaps the row order of "Montserrat" and "Islas Malvinas",
# and bumps the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated statistics for several countries (columns: B,C,D,E,F,G,H) ---

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 5646995
$ws.Cells.Item(4, 3).Value = 35020
$ws.Cells.Item(4, 4).Value = 2999436
$ws.Cells.Item(4, 5).Value = 2472790
$ws.Cells.Item(4, 7).Value = 1053
$ws.Cells.Item(4, 8).Value = 174769

# Row 5 - Brasil
$ws.Cells.Item(5, 2).Value = 3407354
$ws.Cells.Item(5, 3).Value = 44119
$ws.Cells.Item(5, 5).Value = 818972
$ws.Cells.Item(5, 7).Value = 1234
$ws.Cells.Item(5, 8).Value = 109888

# Row 8 - Sudafrica
$ws.Cells.Item(8, 2).Value = 592144
$ws.Cells.Item(8, 3).Value = 2258
$ws.Cells.Item(8, 4).Value = 485468
$ws.Cells.Item(8, 5).Value = 94412
$ws.Cells.Item(8, 7).Value = 282
$ws.Cells.Item(8, 8).Value = 12264

# Row 22 - Alemania
$ws.Cells.Item(22, 2).Value = 228105
$ws.Cells.Item(22, 3).Value = 1419
$ws.Cells.Item(22, 5).Value = 15900

# Row 29 - Kazajistan
$ws.Cells.Item(29, 5).Value = 17440
$ws.Cells.Item(29, 7).Value = 21
$ws.Cells.Item(29, 8).Value = 1415

# Row 33 - Israel
$ws.Cells.Item(33, 2).Value = 96409
$ws.Cells.Item(33, 3).Value = 1658
$ws.Cells.Item(33, 4).Value = 71990
$ws.Cells.Item(33, 5).Value = 23711
$ws.Cells.Item(33, 7).Value = 16
$ws.Cells.Item(33, 8).Value = 708

# Row 65 - Moldavia
$ws.Cells.Item(65, 4).Value = 21885
$ws.Cells.Item(65, 5).Value = 7996

# Row 75 - Camerun
$ws.Cells.Item(75, 2).Value = 18599
$ws.Cells.Item(75, 3).Value = 17
$ws.Cells.Item(75, 5).Value = 1653
$ws.Cells.Item(75, 7).Value = 3
$ws.Cells.Item(75, 8).Value = 406

# Row 77 - Costa de Marfil
$ws.Cells.Item(77, 2).Value = 17150
$ws.Cells.Item(77, 3).Value = 43
$ws.Cells.Item(77, 4).Value = 14183
$ws.Cells.Item(77, 5).Value = 2857

# Row 110 - Republica de Africa Central
$ws.Cells.Item(110, 2).Value = 4679
$ws.Cells.Item(110, 3).Value = 12
$ws.Cells.Item(110, 4).Value = 1755
$ws.Cells.Item(110, 5).Value = 2863

# Row 120 - Cabo Verde
$ws.Cells.Item(120, 2).Value = 3253
$ws.Cells.Item(120, 3).Value = 50
$ws.Cells.Item(120, 4).Value = 2390
$ws.Cells.Item(120, 5).Value = 827

# Row 127 - Ruanda
$ws.Cells.Item(127, 2).Value = 2577
$ws.Cells.Item(127, 3).Value = 37
$ws.Cells.Item(127, 4).Value = 1683
$ws.Cells.Item(127, 5).Value = 884
$ws.Cells.Item(127, 7).Value = 2
$ws.Cells.Item(127, 8).Value = 10

# Row 138 - Sierra Leona
$ws.Cells.Item(138, 2).Value = 1959
$ws.Cells.Item(138, 3).Value = 3
$ws.Cells.Item(138, 4).Value = 1529
$ws.Cells.Item(138, 5).Value = 361

# --- Swap the order of "Islas Malvinas" and "Montserrat" rows (213/214) ---
# Row 213 used to be Islas Malvinas (D=13,H=0); row 214 used to be Montserrat (D=12,H=1).
# After the edit, Montserrat comes first (row 213) and Islas Malvinas second (row 214),
# each keeping its own D/H values.
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 18 de Agosto de 2020 a las 23:20"
